$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.631.08"
$ws.Range("E2").Value = "  -7.73%  "

$ws.Range("D3").Value = "2.550.61"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.08"
$ws.Range("E5").Value = "  -4.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.34"
$ws.Range("E6").Value = "  -8.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -4.27%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -6.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.70"
$ws.Range("E10").Value = "  -8.75%  "

$ws.Range("E11").Value = "  -4.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -5.86%  "

$ws.Range("D13").Value = "2.937.53"
$ws.Range("E13").Value = "  -2.10%  "

$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "2.527.09"
$ws.Range("E15").Value = "  -2.91%  "

$ws.Range("E16").Value = "  -5.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.08"
$ws.Range("E17").Value = "  -5.50%  "

$ws.Range("D18").Value = "42.739.18"
$ws.Range("E18").Value = "  -7.79%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("E20").Value = "  -2.82%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0965"
$ws.Range("E21").Value = "  -5.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.02"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.45"
$ws.Range("E23").Value = "  -6.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -5.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.44"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("E26").Value = "  -4.59%  "

$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("E28").Value = "  -7.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.86"
$ws.Range("E30").Value = "  -6.18%  "

$ws.Range("E31").Value = "  -6.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  -5.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.18"
$ws.Range("E33").Value = "  -3.69%  "

$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("E36").Value = "  -5.82%  "

$ws.Range("E37").Value = "  -8.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.17"
$ws.Range("E38").Value = "  +3.90%  "

$ws.Range("E39").Value = "  -3.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.84"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -5.96%  "

$ws.Range("E42").Value = "  -7.29%  "

$ws.Range("D43").Value = "2.072.66"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("E44").Value = "  -4.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "84.12"
$ws.Range("E46").Value = "  -11.64%  "

$ws.Range("E47").Value = "  +2.12%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.792.63"
$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.75"
$ws.Range("E49").Value = "  -9.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").Value = "  -3.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.11"
$ws.Range("E51").Value = "  -5.81%  "
